# "fix: Send emails was fixed"
# - Row 4 (idContrato 5050): the contact "Alejandro Sani" / 9909 / (no working
#   mailto link) is replaced by "Italo Pilatasig" whose e-mail
#   (bfabita@hotmail.es) now gets a working mailto: hyperlink.
# - The old mailto: hyperlink on D2 (which pointed at the generic
#   alexandermacash@gmail.com address) is removed - the text stays, just the
#   link is gone.
# - A brand-new row 5 is appended for "David Flores" with his own working
#   mailto: hyperlink.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- remove the stale hyperlink on D2 (text/value is untouched) ---
$ws.Range("D2").Hyperlinks.Delete()

# --- row 4: swap the old contact for Italo Pilatasig ---
$ws.Range("B4").Value = "Italo Pilatasig"
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "9909"
$ws.Range("D4").Value = "bfabita@hotmail.es"
$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:bfabita@hotmail.es") | Out-Null
$ws.Range("D4").Style = "Hipervínculo"

# --- row 5 (new): David Flores ---
$ws.Range("A5").Value = 5050
$ws.Range("B5").Value = "David Flores"
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "984556639"
$ws.Range("D5").Value = "fxbricio7@gmail.com"
$ws.Range("E5").NumberFormat = "$#,##0.00;[Red]$-#,##0.00"
$ws.Range("E5").Value = 5
$ws.Range("F5").Value = "GEDATECU SA"
$ws.Hyperlinks.Add($ws.Range("D5"), "mailto:fxbricio7@gmail.com") | Out-Null
$ws.Range("D5").Style = "Hipervínculo"

# --- match the saved selection/view from the commit ---
$ws.Range("E6").Select() | Out-Null
